$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Transaction details to download/populate (Category, Amount, Date-serial)
$rows = @(
    @("Drinks",         2765, 45991.22928240741),
    @("Fruits",          135, 45989.22928240741),
    @("Light Bill",      586, 45985.22928240741),
    @("Drinks",         1560, 45985.22928240741),
    @("Birthday Party", 2056, 45982.22928240741),
    @("Dinner",          347, 45979.22928240741),
    @("Fuel",           1100, 45979.22928240741),
    @("Rent",           2400, 45974.22928240741),
    @("Bike Servicing", 1895, 45973.22928240741),
    @("Dinner",          462, 45970.22928240741)
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}

# Propagate the existing date format from C2 down through the new rows
$ws.Range("C2").Copy()
$ws.Range("C2:C11").PasteSpecial(-4122)
